$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows in each "Category" group (Level of education, Affiliation type,
# Position) were reordered. Update each affected row's Option label and
# counts to reflect the new row order; numeric/label values for rows not
# listed here are unchanged.

# Level of education group
$ws.Range("B3").Value = "Other (please specify below)"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.2

$ws.Range("B5").Value = "Bachelor’s degree"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 12.9

# Affiliation type group
$ws.Range("B6").Value = "Government"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.2

$ws.Range("B7").Value = "Industry"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 6.5

# Position group
$ws.Range("B9").Value = "Graduate student (including professional school student)"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = 12.9

$ws.Range("B10").Value = "Undergraduate student"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.2

$ws.Range("B11").Value = "Postdoc"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 12.9

$ws.Range("B12").Value = "Faculty member"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 10
$ws.Range("G12").Value = 32.3

$ws.Range("B14").Value = "Staff member (including research/academic/teaching staff)"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 10
$ws.Range("G14").Value = 32.3
